$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.265.21'
$ws.Range("E2").Value = '  -5.66%  '

$ws.Range("D3").Value = '1.835.60'
$ws.Range("E3").Value = '  -5.41%  '

$ws.Range("E4").Value = '  -0.79%  '

$ws.Range("D5").Value = "'331.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.14%  '

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.79%  '

$ws.Range("D7").Value = "'0.4600"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.06%  '

$ws.Range("D8").Value = "'0.3860"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.79%  '

$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = "'45.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.99%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = "'0.07843"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.48%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = "'0.9647"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.14%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = "'21.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -8.48%  '

$ws.Range("D13").Value = "'5.721"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.28%  '

$ws.Range("D14").Value = "'6.913"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.51%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.757.88'
$ws.Range("E15").Value = '  -10.08%  '

$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").Value = "'0.06879"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = "'1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.86%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = "'86.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.26%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = "'0.000009921"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.64%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = "'16.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.66%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.73%  '

$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").Value = '28.274.28'
$ws.Range("E22").Value = '  -5.63%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = "'5.333"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.58%  '

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = "'10.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.03%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = "'2.156"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.78%  '

$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.053.61'
$ws.Range("E26").Value = '  -6.23%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = "'153.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.58%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'19.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.68%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = "'5.812"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -11.93%  '

$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").Value = "'1.974"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.26%  '

$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = "'116.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.85%  '

$ws.Range("D32").Value = "'0.9427"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.31%  '

$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = "'0.09312"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.40%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = "'5.284"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.07%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = "'3.443"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.31%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = "'1.325"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.81%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = "'0.06034"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.98%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = "'0.02158"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.96%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = "'1.154"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.54%  '

$ws.Range("B40").Value = 'Frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D40").Value = "'1.000"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.78%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = "'0.5615"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.29%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'7.552"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.92%  '

$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = "'9.995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.03%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = "'0.1781"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.09%  '

$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = "'1.241"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.99%  '

$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = "'2.269"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.63%  '

$ws.Range("D47").Value = "'11.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.05%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = "'0.5297"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.05%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = "'0.07026"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.64%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = "'1.838"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.48%  '

$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = "'113.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.11%  '
